$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text storage for price cells that would otherwise parse as numbers
$priceTextCells = @('D5', 'D6', 'D9', 'D10', 'D11', 'D12', 'D13', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D38', 'D39', 'D42', 'D43', 'D45', 'D47', 'D48', 'D49', 'D51')
foreach ($addr in $priceTextCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '67.941.15'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '3.649.35'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '592.79'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').Value = '189.85'
$ws.Range('E6').Value = '  +4.64%  '
$ws.Range('E7').Value = '  -1.51%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '0.694'
$ws.Range('E9').Value = '  -3.22%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.151'
$ws.Range('E10').Value = '  -7.43%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').Value = '56.31'
$ws.Range('E11').Value = '  +2.73%  '
$ws.Range('D12').Value = '0.0000268'
$ws.Range('E12').Value = '  -8.05%  '
$ws.Range('D13').Value = '10.10'
$ws.Range('E13').Value = '  -2.97%  '
$ws.Range('D14').Value = '4.242.39'
$ws.Range('E14').Value = '  +0.73%  '
$ws.Range('D15').Value = '3.656.52'
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').Value = '18.75'
$ws.Range('E17').Value = '  -3.79%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '67.779.15'
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').Value = '1.10'
$ws.Range('E19').Value = '  -2.52%  '
$ws.Range('D20').Value = '12.34'
$ws.Range('E20').Value = '  -3.29%  '
$ws.Range('D21').Value = '397.49'
$ws.Range('E21').Value = '  -2.83%  '
$ws.Range('D22').Value = '4.35'
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('D23').Value = '87.27'
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('E24').Value = '  -3.07%  '
$ws.Range('D25').Value = '10.89'
$ws.Range('E25').Value = '  -2.69%  '
$ws.Range('D26').Value = '12.42'
$ws.Range('E26').Value = '  -3.02%  '
$ws.Range('D27').Value = '6.06'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = '3.63'
$ws.Range('E28').Value = '  -7.42%  '
$ws.Range('D29').Value = '9.25'
$ws.Range('E29').Value = '  -2.38%  '
$ws.Range('D30').Value = '31.60'
$ws.Range('E30').Value = '  -3.13%  '
$ws.Range('D31').Value = '7.21'
$ws.Range('E31').Value = '  -3.14%  '
$ws.Range('D32').Value = '12.21'
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('D33').Value = '65.71'
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('D34').Value = '43.86'
$ws.Range('E34').Value = '  +0.87%  '
$ws.Range('D35').Value = '603.15'
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.13%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '0.391'
$ws.Range('E39').Value = '  -2.56%  '
$ws.Range('D40').Value = '0.0₃0756'
$ws.Range('E40').Value = '  -15.79%  '
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('D42').Value = '2.87'
$ws.Range('E42').Value = '  -4.60%  '
$ws.Range('D43').Value = '0.0423'
$ws.Range('E43').Value = '  -2.98%  '
$ws.Range('E44').Value = '  -9.48%  '
$ws.Range('D45').Value = '0.134'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('D46').Value = '2.771.51'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').Value = '3.11'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').Value = '143.16'
$ws.Range('E48').Value = '  +2.74%  '
$ws.Range('D49').Value = '8.67'
$ws.Range('E49').Value = '  -6.72%  '
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('D51').Value = '2.45'
$ws.Range('E51').Value = '  -16.44%  '

$wb.Save()
